# Generate Report for Handoff
# Replace the old GUID-based file names/content hashes with the new ones,
# and bump the associated timestamps, across the three worksheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "54c1bc11-4892-4edd-81d5-a632eb066f2b"
$newGuid = "4509f7f9-3db0-4973-aa2d-e49ac65dd39d"
$oldHash = "2127df2aa1eb8a2cc2aff140002115c636436d32"
$newHash = "930326b3155782b12a5fbd62c6687155637795f2"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The hyperlink target (Address) itself is unchanged by the commit - only the
# visible text of each hyperlink (and the backing cell value) changes to
# reflect the new GUID-based file name. This engine's Hyperlink objects only
# support clean updates via delete-then-recreate (in-place property writes on
# an existing Hyperlink duplicate the link), so rebuild each one explicitly.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a62beaf1062d44b94d5e92f114d8b79138f4123/e2e/$oldGuid.md"

# --- Overview sheet ---
# A2: plain markdown file name
$wsOverview.Range("A2").Value = "$newGuid.md"
# B2: hyperlink cell - update the cell text and recreate the hyperlink so its
# displayed text matches, keeping the same link target as before.
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", "e2e\$newGuid.md")
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-21 11:01:52"

# --- zh-cn sheet ---
$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 11:01:48"

# --- de-de sheet ---
$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 11:01:52"
